$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("F2").Value = 2.54

# Row 3 changes
$ws.Range("H3").Value = 10
$ws.Range("Q3").Value = 2.1
$ws.Range("T3").Value = 2.42
$ws.Range("U3").Value = 1.66
$ws.Range("AA3").Value = 590
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 40
$ws.Range("AI3").Value = 230
$ws.Range("AM3").Value = 320

# Row 4 changes
$ws.Range("F4").Value = 1.61
$ws.Range("G4").Value = 1.62
$ws.Range("H4").Value = 6.2
$ws.Range("I4").Value = 6.4
$ws.Range("J4").Value = 4.5
$ws.Range("K4").Value = 4.6
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 1.81
$ws.Range("R4").Value = 1.46
$ws.Range("T4").Value = 1.88
$ws.Range("X4").Value = 19
$ws.Range("Y4").Value = 22
$ws.Range("Z4").Value = 55
$ws.Range("AA4").Value = 210
$ws.Range("AB4").Value = 9
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 24
$ws.Range("AE4").Value = 85
$ws.Range("AF4").Value = 9.6
$ws.Range("AI4").Value = 75
$ws.Range("AJ4").Value = 15.5
$ws.Range("AN4").Value = 8
$ws.Range("AO4").Value = 120
